$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")  # Sheet1 is also $wb.ActiveSheet in this workbook

# New LeetCode entry: day 31 - Binary Tree Preorder Traversal
$ws.Range("A33").Value = 31

# Add the hyperlink for the problem name (creates the shared-string + relationship),
# matching the pattern used by every other row in this column.
$ws.Hyperlinks.Add($ws.Range("B33"), "https://leetcode.com/problems/binary-tree-preorder-traversal/", [Type]::Missing, [Type]::Missing, "https://leetcode.com/problems/binary-tree-preorder-traversal/") | Out-Null
# Restore the visible cell text to the problem name and reapply the Hyperlink cell style
# (Hyperlinks.Add() both overwrites the display text and re-derives a fresh style).
$ws.Range("B33").Value = "Binary Tree Preorder Traversal"
$ws.Range("B33").Style = "Hyperlink"

$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 28
$ws.Range("F33").Value = 0.0246
$ws.Range("G33").Value = 16.23
$ws.Range("H33").Value = 0.43
$ws.Range("I33").Value = "https://leetcode.com/problems/binary-tree-preorder-traversal/submissions/"

# Move the active selection to I33, matching where the author finished editing.
$ws.Range("I33").Select()
